# Inicio lectura de coordinadora
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7: manuscript status text changed
$ws.Range("F7").Value = "Manuscritos en revisión nuevo editor"

# Row 8: underlying shared string text changed (same pointer, new wording)
$ws.Range("F8").Value = "Manuscritos en revisión nuevo editor"

# Row 17: new status text (entered before rows 12/13/16 per editing order)
$ws.Range("F17").Value = "En ajustes pedidos por coordinadora"

# Row 16: date + new status text
$ws.Range("E16").Value = [DateTime]"2015-03-16"
$ws.Range("F16").Value = "En digitación"

# Row 12: dates + new status text
$ws.Range("B12").Value = [DateTime]"2015-03-17"
$ws.Range("C12").Value = [DateTime]"2015-03-17"
$ws.Range("D12").Value = [DateTime]"2015-03-17"
$ws.Range("E12").Value = [DateTime]"2015-03-17"
$ws.Range("F12").Value = "En búsqueda gráfica"

# Row 13: new status text (last edit, leaves selection here)
$ws.Range("F13").Value = "En manuscrito"

# Row heights re-wrap to fit the new (shorter/longer) text
$ws.Rows.Item(8).RowHeight = 30.75
$ws.Rows.Item(17).RowHeight = 30.75

# Update the active selection to reflect where editing ended
$ws.Range("F13").Select() | Out-Null
